$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New header for column M
$ws.Range("M1").Value = "MatchOutcome"

# Format the outcome column as Text before writing values so "1"/"2" stay
# text instead of being coerced to numbers.
$ws.Range("M2:M39").NumberFormat = "@"

# Seed the shared-string table in the same order the source workbook used
# (X, then 1, then 2) by writing the first occurrence of each value first.
$ws.Range("M6").Value = "X"
$ws.Range("M3").Value = "1"
$ws.Range("M2").Value = "2"

# Remaining match-outcome values.
$ws.Range("M4").Value = "1"
$ws.Range("M5").Value = "2"
$ws.Range("M10").Value = "1"
$ws.Range("M13").Value = "1"
$ws.Range("M14").Value = "X"
$ws.Range("M18").Value = "X"

# Widen column M to fit the new header/values (closest the host's pixel
# rounding allows to the source width of 25.7109375).
$ws.Columns.Item(13).ColumnWidth = 24.76

# Match the author's final selection/view state.
$ws.Range("M8").Select()
